$d = $word.ActiveDocument
$d.Content.Find.Execute("September 20, 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "October 14, 2024", 2)
